# "Generate Report for Handback"
#
# The localization file ec47e00c-272c-40d5-b36e-3c243d969b9d has come back
# from handback (in sync with en-US) for both the zh-cn and de-de locales.
# Update the Overview sheet's status for that file, and fill in each locale
# sheet's "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns (F, G, H) for that file's row.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: file ec47e00c...md is row 2 (zh-cn column B, de-de column C)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: row 2 is the ec47e00c...md file
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusHandedBack

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/677beb26ad75a95c78894b707cac7839f52d48e8/e2e/ec47e00c-272c-40d5-b36e-3c243d969b9d.md",
    "",
    "",
    "ec47e00c-272c-40d5-b36e-3c243d969b9d.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/645c75ebc5626ee8cee4c5f0a327b8e4f774060e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ec47e00c-272c-40d5-b36e-3c243d969b9d.088e788414899d7202b0236b9fb06bdcb8cddbb8.zh-cn.xlf",
    "",
    "",
    "ec47e00c-272c-40d5-b36e-3c243d969b9d.088e788414899d7202b0236b9fb06bdcb8cddbb8.zh-cn.xlf"
) | Out-Null

$zhcn.Range("H2").Value = "2016-03-17 22:29:39"

# ---------------------------------------------------------------------
# de-de sheet: row 2 is the ec47e00c...md file
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusHandedBack

$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/677beb26ad75a95c78894b707cac7839f52d48e8/e2e/ec47e00c-272c-40d5-b36e-3c243d969b9d.md",
    "",
    "",
    "ec47e00c-272c-40d5-b36e-3c243d969b9d.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35c68a050e110761da2a3af41ccaa17c5f1fd32b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ec47e00c-272c-40d5-b36e-3c243d969b9d.088e788414899d7202b0236b9fb06bdcb8cddbb8.de-de.xlf",
    "",
    "",
    "ec47e00c-272c-40d5-b36e-3c243d969b9d.088e788414899d7202b0236b9fb06bdcb8cddbb8.de-de.xlf"
) | Out-Null

$dede.Range("H2").Value = "2016-03-17 22:29:45"
